# edit.ps1 - applies the two substantive changes from the commit:
#
#  1. Slide 6's table (shape 2) switches its table style (tableStyleId)
#     from {232923E1-7C02-4AFD-8B70-6CA34679894B} to
#     {CD09AD97-8229-49C1-9173-A0DA38CFEE24}.
#
#  2. The theme that the slide master actually renders with (the part
#     physically named ppt/theme/theme2.xml, which currently carries the
#     "Integral" colour palette) is switched over to the "Office Theme"
#     colour palette (the palette that currently lives in
#     ppt/theme/theme1.xml, which is only ever used by the notes master).
#     This is done by rewriting each of the 12 theme colour slots on the
#     slide master's theme colour scheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap on slide 6.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{CD09AD97-8229-49C1-9173-A0DA38CFEE24}")

# ---------------------------------------------------------------------
# 2) Slide master theme colour swap (Integral -> Office Theme values).
# ---------------------------------------------------------------------
function ConvertTo-VbaRgb($hex) {
    # VBA/COM RGB() integer packs as R + G*256 + B*65536 (matches the
    # RGB property exposed on ThemeColorScheme color slots).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the 12 ThemeColorScheme slots: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-VbaRgb $officeThemeColors[$i - 1]
}
